$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column / row header labels from underscore style to hyphen style
$ws.Range("B1").Value = "Col-01"
$ws.Range("C1").Value = "Col-02"
$ws.Range("D1").Value = "Col-03"
$ws.Range("E1").Value = "Col-04"

$ws.Range("A2").Value = "Row-01"
$ws.Range("A3").Value = "Row-02"
$ws.Range("A4").Value = "Row-03"
$ws.Range("A5").Value = "Row-04"

# Update the saved selection / active cell shown in the sheet view
$ws.Range("E10").Select()
